$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.035.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.48%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.420.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'552.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.17%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'137.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.43%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +2.81%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +2.33%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'5.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.60%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.64%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.79%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'24.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +4.76%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.851.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.29%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'59.974.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.53%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +2.42%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.408.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.54%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'11.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +5.64%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +1.89%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'331.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.58%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.24%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D23").Value = "'65.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.89%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +3.49%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'8.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.94%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.20%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.99%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0784"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +6.67%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.60%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'170.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.17%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.34%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'18.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.53%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +2.13%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.01%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +5.38%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D37").Value = "'4.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.21%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.45%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +1.29%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.415"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +9.72%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'314.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +8.25%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.21%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'140.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.21%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0964"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.35%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0520"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.77%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "Polygon"
$ws.Range("C46").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D46").Value = "'0.412"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +7.98%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'19.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.87%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.576"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.77%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0226"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.52%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'17.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.25%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'11.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.37%  "
$ws.Range("E51").Style = "Normal"
